# Atualizado por script em 03-11-2023 14:45
# Append the new match row (row 70) to the Armenia Premier League sheet,
# mirroring the formatting of the previous data row (row 69).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 70
$prevRow = $newRow - 1

# Copy the formatting (styles) of the previous row down onto the new row
# so the new cells reuse the existing style indices (bold/border index cell,
# datetime-formatted match-date cell, etc.) instead of creating new ones.
$ws.Range("A" + $prevRow + ":V" + $prevRow).Copy()
$ws.Range("A" + $newRow + ":V" + $newRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new match data.
$ws.Cells.Item($newRow, 1).Value2 = 69
$ws.Cells.Item($newRow, 2).Value2 = "armenia"
$ws.Cells.Item($newRow, 3).Value2 = "premier-league"
$ws.Cells.Item($newRow, 4).Value2 = "2023-2024"
$ws.Cells.Item($newRow, 5).Value2 = 45233.54166666666
$ws.Cells.Item($newRow, 6).Value2 = "Urartu"
$ws.Cells.Item($newRow, 7).Value2 = 1
$ws.Cells.Item($newRow, 8).Value2 = "West Armenia"
$ws.Cells.Item($newRow, 9).Value2 = 2
$ws.Cells.Item($newRow, 10).Value2 = 1.02
$ws.Cells.Item($newRow, 11).Value2 = "03/11/2023 11:27"
$ws.Cells.Item($newRow, 12).Value2 = 1.02
$ws.Cells.Item($newRow, 13).Value2 = "03/11/2023 11:27"
$ws.Cells.Item($newRow, 14).Value2 = 17.67
$ws.Cells.Item($newRow, 15).Value2 = "03/11/2023 11:30"
$ws.Cells.Item($newRow, 16).Value2 = 17.67
$ws.Cells.Item($newRow, 17).Value2 = "03/11/2023 11:30"
$ws.Cells.Item($newRow, 18).Value2 = 27.78
$ws.Cells.Item($newRow, 19).Value2 = "03/11/2023 11:30"
$ws.Cells.Item($newRow, 20).Value2 = 27.78
$ws.Cells.Item($newRow, 21).Value2 = "03/11/2023 11:30"
$ws.Cells.Item($newRow, 22).Value2 = "https://www.betexplorer.com/football/armenia/premier-league/urartu-west-armenia/xdYDlZd6/"
